$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.198.21"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "2.416.68"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.531"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.76%  "
$ws.Range("D9").Value = "2.410.30"
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("E11").Value = "  -2.24%  "
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").Value = "2.851.99"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").Value = "62.058.74"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "2.412.07"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "323.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.41%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.81%  "
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "576.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.21%  "
$ws.Range("E28").Value = "  +3.40%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("B30").Value = "WrappedeETH"
$ws.Range("C30").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D30").Value = "2.531.16"
$ws.Range("E30").Value = "  +2.32%  "
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  -3.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.382"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "152.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.61%  "
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("E47").Value = "  +1.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("E49").Value = "  +1.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0918"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("E51").Value = "  +1.64%  "
